$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.00" or
# "577.61" are not auto-converted to numbers by Excel type inference, matching
# the original inline-string cell type. Wrapping the whole used range keeps
# this to a single style allocation instead of one per edited cell.
$dPrices = $ws.Range("D2:D51")
$dPrices.NumberFormat = "@"

$ws.Range("D2").Value = '65.214.29'
$ws.Range("E2").Value = '  +4.06%  '
$ws.Range("D3").Value = '3.483.34'
$ws.Range("E3").Value = '  +3.67%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '577.61'
$ws.Range("E5").Value = '  +3.36%  '
$ws.Range("D6").Value = '161.46'
$ws.Range("E6").Value = '  +5.13%  '
$ws.Range("D7").Value = '0.616'
$ws.Range("E7").Value = '  +15.01%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '3.483.75'
$ws.Range("E9").Value = '  +3.66%  '
$ws.Range("D10").Value = '7.26'
$ws.Range("E10").Value = '  -1.17%  '
$ws.Range("D11").Value = '0.125'
$ws.Range("E11").Value = '  +4.62%  '
$ws.Range("D12").Value = '0.446'
$ws.Range("E12").Value = '  +4.00%  '
$ws.Range("D13").Value = '4.095.37'
$ws.Range("E13").Value = '  +3.88%  '
$ws.Range("E14").Value = '  +1.04%  '
$ws.Range("D15").Value = '0.0000193'
$ws.Range("E15").Value = '  +3.75%  '
$ws.Range("D16").Value = '28.68'
$ws.Range("E16").Value = '  +7.47%  '
$ws.Range("D17").Value = '65.339.91'
$ws.Range("E17").Value = '  +4.13%  '
$ws.Range("D18").Value = '3.459.22'
$ws.Range("E18").Value = '  +5.02%  '
$ws.Range("D19").Value = '6.46'
$ws.Range("E19").Value = '  +4.64%  '
$ws.Range("D20").Value = '14.33'
$ws.Range("E20").Value = '  +3.04%  '
$ws.Range("D21").Value = '382.10'
$ws.Range("E21").Value = '  +2.67%  '
$ws.Range("D22").Value = '8.21'
$ws.Range("E22").Value = '  +3.45%  '
$ws.Range("D23").Value = '0.550'
$ws.Range("E23").Value = '  +5.30%  '
$ws.Range("D24").Value = '72.88'
$ws.Range("E24").Value = '  +3.03%  '
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("D26").Value = '0.0000119'
$ws.Range("E26").Value = '  +5.50%  '
$ws.Range("D27").Value = '10.06'
$ws.Range("E27").Value = '  +7.22%  '
$ws.Range("E28").Value = '  +2.73%  '
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.40%  '
$ws.Range("D30").Value = '1.52'
$ws.Range("E30").Value = '  +14.18%  '
$ws.Range("D31").Value = '6.24'
$ws.Range("E31").Value = '  +3.90%  '
$ws.Range("D32").Value = '2.05'
$ws.Range("E32").Value = '  +4.67%  '
$ws.Range("D33").Value = '23.59'
$ws.Range("E33").Value = '  +2.83%  '
$ws.Range("D34").Value = '7.24'
$ws.Range("E34").Value = '  +8.50%  '
$ws.Range("D35").Value = '1.58'
$ws.Range("E35").Value = '  +10.30%  '
$ws.Range("D36").Value = '161.59'
$ws.Range("E36").Value = '  +1.67%  '
$ws.Range("D37").Value = '1.92'
$ws.Range("E37").Value = '  +6.63%  '
$ws.Range("D38").Value = '3.033.04'
$ws.Range("E38").Value = '  +4.77%  '
$ws.Range("D39").Value = '0.0775'
$ws.Range("E39").Value = '  +2.03%  '
$ws.Range("D40").Value = '26.98'
$ws.Range("E40").Value = '  +0.64%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '4.57'
$ws.Range("E41").Value = '  +7.00%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.0322'
$ws.Range("E42").Value = '  +2.47%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '6.52'
$ws.Range("E43").Value = '  +0.77%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = '42.84'
$ws.Range("E44").Value = '  +4.48%  '
$ws.Range("D45").Value = '0.779'
$ws.Range("E45").Value = '  +5.44%  '
$ws.Range("D46").Value = '25.90'
$ws.Range("E46").Value = '  +13.84%  '
$ws.Range("D47").Value = '1.11'
$ws.Range("E47").Value = '  +5.79%  '
$ws.Range("D48").Value = '319.33'
$ws.Range("E48").Value = '  +12.68%  '
$ws.Range("D49").Value = '6.76'
$ws.Range("E49").Value = '  +7.22%  '
$ws.Range("E50").Value = '  +8.05%  '
$ws.Range("D51").Value = '2.20'
$ws.Range("E51").Value = '  +5.77%  '

# Restore the default style on column D so every cell keeps the workbook
# original (style-less) formatting.
$dPrices.Style = "Normal"
